# Re-save round-trip (this workbook was re-opened/re-saved from a
# different machine/Excel build) resulted in two visible, intentional
# changes that this script reproduces:
#
#   1. The second worksheet's truncated tab name "Facility Fundamen..."
#      is corrected to its full name "Facility Fundamentals".
#   2. The active worksheet/selection moves off "Test Survey" (E1) onto
#      "Facility Fundamentals", with its selection moved from G7 to E21.

$wb = $excel.ActiveWorkbook

# 1. Fix the truncated worksheet name (fall back to positional index 2
#    if the name has already been changed / doesn't match exactly).
$wsFacility = $null
foreach ($s in $wb.Worksheets) {
    if ($s.Name -eq "Facility Fundamen...") {
        $wsFacility = $s
    }
}
if ($wsFacility -eq $null) {
    $wsFacility = $wb.Worksheets.Item(2)
}
$wsFacility.Name = "Facility Fundamentals"

# 2. Make "Facility Fundamentals" the active sheet and move its selection.
$wsFacility.Activate()
$wsFacility.Range("E21").Select()
